$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A2:Z13").Select() | Out-Null
$ws.Columns("AA").Delete()
$ws.Columns("U").Delete()
$ws.Columns("M").Delete()
$ws.Columns("L").Delete()
$ws.Columns("F").Delete()
